$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = 3.3
$ws.Range("J2").Value = 2.6
$ws.Range("L2").Value = 3.75
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.38
$ws.Range("S2").Value = 1.99
$ws.Range("T2").Value = 1.91
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.57
$ws.Range("W2").Value = 1.29
$ws.Range("X2").Value = 3.5
$ws.Range("AC2").Value = 9.5
$ws.Range("AE2").Value = 15

# Row 3
$ws.Range("G3").Value = 2.4
$ws.Range("I3").Value = 2.8
$ws.Range("Q3").Value = 1.9
$ws.Range("R3").Value = 2
$ws.Range("AA3").Value = 10
$ws.Range("AL3").Value = 11
$ws.Range("AM3").Value = 17
$ws.Range("AN3").Value = 12

# Row 4
$ws.Range("G4").Value = 2.63
$ws.Range("I4").Value = 2.6
$ws.Range("J4").Value = 3.25
$ws.Range("L4").Value = 3.25
$ws.Range("AL4").Value = 10

# Row 5
$ws.Range("G5").Value = 1.18
$ws.Range("H5").Value = 8
$ws.Range("J5").Value = 1.5
$ws.Range("L5").Value = 9.5
$ws.Range("Q5").Value = 1.29
$ws.Range("R5").Value = 3.75
$ws.Range("Y5").Value = 1.8
$ws.Range("Z5").Value = 1.95
$ws.Range("AD5").Value = 9
$ws.Range("AF5").Value = 26
$ws.Range("AI5").Value = 26
$ws.Range("AK5").Value = 251
$ws.Range("AN5").Value = 41
$ws.Range("AO5").Value = 201
$ws.Range("AQ5").Value = 67

# Row 7
$ws.Range("G7").Value = 1.67
$ws.Range("H7").Value = 3.3
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 2.4
$ws.Range("L7").Value = 6.5
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5
$ws.Range("Q7").Value = 2.5
$ws.Range("R7").Value = 1.5
$ws.Range("S7").Value = 4
$ws.Range("Y7").Value = 2.38
$ws.Range("Z7").Value = 1.53
$ws.Range("AB7").Value = 6.5
$ws.Range("AD7").Value = 12
$ws.Range("AE7").Value = 17
$ws.Range("AH7").Value = 7
$ws.Range("AI7").Value = 23
$ws.Range("AJ7").Value = 101
$ws.Range("AL7").Value = 11
$ws.Range("AM7").Value = 29
$ws.Range("AN7").Value = 21
$ws.Range("AO7").Value = 67
$ws.Range("AQ7").Value = 67
$ws.Range("AR7").Value = 1.93
$ws.Range("AS7").Value = 1.93

# Row 8
$ws.Range("N8").Value = 9

# Row 10
$ws.Range("G10").Value = 1.65
$ws.Range("I10").Value = 5.25
$ws.Range("J10").Value = 2.3
$ws.Range("L10").Value = 5.5
$ws.Range("Q10").Value = 2.05
$ws.Range("R10").Value = 1.75
$ws.Range("Y10").Value = 2
$ws.Range("Z10").Value = 1.75
$ws.Range("AK10").Value = 401
$ws.Range("AL10").Value = 12
$ws.Range("AM10").Value = 26
$ws.Range("AN10").Value = 17
$ws.Range("AO10").Value = 51

# Row 11
$ws.Range("I11").Value = 2.05
$ws.Range("AH11").Value = 6.5

# Row 12
$ws.Range("G12").Value = 2.88
$ws.Range("I12").Value = 2.5
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.5
$ws.Range("Q12").Value = 2.5
$ws.Range("R12").Value = 1.5
$ws.Range("Y12").Value = 2.1
$ws.Range("Z12").Value = 1.67
$ws.Range("AB12").Value = 13
$ws.Range("AE12").Value = 29
$ws.Range("AR12").Value = 1.93
$ws.Range("AS12").Value = 1.93

# Row 13
$ws.Range("G13").Value = 1.8
$ws.Range("I13").Value = 3.7
$ws.Range("AB13").Value = 10
$ws.Range("AK13").Value = 151
$ws.Range("AP13").Value = 29

# Row 15
$ws.Range("G15").Value = 4.1
$ws.Range("H15").Value = 3.9
$ws.Range("I15").Value = 1.73
$ws.Range("J15").Value = 4.5
$ws.Range("L15").Value = 2.25
$ws.Range("M15").Value = 1.01
$ws.Range("N15").Value = 13
$ws.Range("AB15").Value = 23
$ws.Range("AC15").Value = 15
$ws.Range("AE15").Value = 34
$ws.Range("AL15").Value = 8
$ws.Range("AO15").Value = 13

# Row 16
$ws.Range("AG16").Value = 29
$ws.Range("AJ16").Value = 51
$ws.Range("AL16").Value = 41
$ws.Range("AM16").Value = 67
$ws.Range("AN16").Value = 34

# Row 17
$ws.Range("M17").Value = 19
$ws.Range("N17").Value = 1.03
$ws.Range("W17").Value = 1.25
$ws.Range("X17").Value = 3.75
$ws.Range("AK17").Value = 151
$ws.Range("AL17").Value = 9.5
$ws.Range("AM17").Value = 8.5
$ws.Range("AO17").Value = 11

# Row 18
$ws.Range("G18").Value = 2.2
$ws.Range("I18").Value = 3
$ws.Range("J18").Value = 2.75
$ws.Range("L18").Value = 3.5
$ws.Range("M18").Value = 1.02
$ws.Range("N18").Value = 12
$ws.Range("AB18").Value = 11
$ws.Range("AO18").Value = 34

# Row 19
$ws.Range("H19").Value = 8.5
$ws.Range("K19").Value = 3.2
$ws.Range("Q19").Value = 1.3
$ws.Range("R19").Value = 3.5
$ws.Range("Y19").Value = 1.95
$ws.Range("Z19").Value = 1.8
$ws.Range("AD19").Value = 7.5
$ws.Range("AK19").Value = 301
$ws.Range("AM19").Value = 67
$ws.Range("AN19").Value = 34
$ws.Range("AP19").Value = 81

# Row 21
$ws.Range("M21").Value = 1.06
$ws.Range("N21").Value = 8
$ws.Range("O21").Value = 1.3
$ws.Range("P21").Value = 3.4
$ws.Range("Q21").Value = 2.05
$ws.Range("R21").Value = 1.75
$ws.Range("U21").Value = 3.5
$ws.Range("V21").Value = 1.29
